$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Preserve formatting for the row that is about to move (old row 8,
#    "order_id") onto its new destination (row 11) BEFORE anything overwrites
#    row 8. Copy cell-by-cell (not as a single A:K range) so that columns
#    which have no cell in the source row (C, F, I, K) don't materialise as
#    blank cells in the destination.
# ---------------------------------------------------------------------------
$ws.Range("A8").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("J8").Copy()
$ws.Range("J11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Apply the fully-populated "items" row style (all 11 columns, A:K) from
#    the existing row 7 ("items"/"price") across the new block of rows
#    7:10 - this gives the new effdate/code/price/discounted_price rows the
#    same wrap/valign formatting as the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A7:K7").Copy()
$ws.Range("A7:K10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Build the new, stand-alone "product_id" row (row 12) by copying the
#    per-cell formatting from row 6 (same sparse A,B,E,G,H,J pattern).
# ---------------------------------------------------------------------------
$ws.Range("A6").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("E6").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("H6").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("J6").Copy()
$ws.Range("J12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Now populate the cell values. Existing table vocabulary (orders_src,
#    orders_tgt, items, string, array, decimal, ...) is reused automatically;
#    brand-new terms are introduced in this order.
# ---------------------------------------------------------------------------

# existing row 4 (customer_id source) - target sub column renamed
$ws.Range("I4").Value = "cust_id"

# new row 12 - product_id passthrough
$ws.Range("A12").Value = "orders_src"
$ws.Range("B12").Value = "product_id"
$ws.Range("E12").Value = "integer"
$ws.Range("G12").Value = "orders_tgt"
$ws.Range("H12").Value = "product_id"
$ws.Range("J12").Value = "integer"

# new row 8 - items/code transformation
$ws.Range("A8").Value = "orders_src"
$ws.Range("B8").Value = "items"
$ws.Range("C8").Value = "code"
$ws.Range("D8").Value = "case when effdate<'2023-10-01' then ""old"" else new"
$ws.Range("E8").Value = "array"
$ws.Range("F8").Value = "string"
$ws.Range("G8").Value = "orders_tgt"
$ws.Range("H8").Value = "item_code"
$ws.Range("J8").Value = "string"

# new row 7 - items/effdate
$ws.Range("A7").Value = "orders_src"
$ws.Range("B7").Value = "items"
$ws.Range("C7").Value = "effdate"
$ws.Range("E7").Value = "array"
$ws.Range("F7").Value = "date"
$ws.Range("G7").Value = "orders_tgt"
$ws.Range("H7").Value = "item_effdate"
$ws.Range("J7").Value = "date"

# new row 10 - items/discounted_price
$ws.Range("A10").Value = "orders_src"
$ws.Range("B10").Value = "items"
$ws.Range("E10").Value = "array"
$ws.Range("F10").Value = "decimal"
$ws.Range("G10").Value = "orders_tgt"
$ws.Range("C10").Value = "discounted_price"
$ws.Range("H10").Value = "item_discount_price"
$ws.Range("D10").Value = "discount price by 15%"

# row 11 - order_id, moved down (same content as the old row 8)
$ws.Range("A11").Value = "orders_src"
$ws.Range("B11").Value = "order_id"
$ws.Range("D11").Value = "concat ""OO5"" before order id"
$ws.Range("E11").Value = "string"
$ws.Range("G11").Value = "orders_tgt"
$ws.Range("H11").Value = "order_id"
$ws.Range("J11").Value = "string"

# row 9 - items/price, moved down (same content as the old row 7)
$ws.Range("A9").Value = "orders_src"
$ws.Range("B9").Value = "items"
$ws.Range("C9").Value = "price"
$ws.Range("E9").Value = "array"
$ws.Range("F9").Value = "decimal"
$ws.Range("G9").Value = "orders_tgt"
$ws.Range("H9").Value = "item_price"
$ws.Range("J9").Value = "decimal"

# ---------------------------------------------------------------------------
# 5) Row heights, to match the auto-fit wrapping of the new text.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 29.5
$ws.Rows.Item(8).RowHeight = 103.25
$ws.Rows.Item(9).RowHeight = 29.5
$ws.Rows.Item(10).RowHeight = 44.25
$ws.Rows.Item(11).RowHeight = 29.5
$ws.Rows.Item(12).RowHeight = 29.5

# ---------------------------------------------------------------------------
# 6) View state: scroll so row 6 is at the top and select D10, matching the
#    saved workbook view.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
